# Increase the outline/indent level by one for four bullet paragraphs in the
# "Area of Interest" content placeholder on slide 1:
#   - "Data architectures and methodologies for Digital Twins"  (lvl 1 -> 2)
#   - "Precision Agriculture:"                                   (lvl 1 -> 2)
#   - "Irrigation optimization"                                  (lvl 2 -> 3)
#   - "Data Platform for italian agriculture domain @Agritech"   (lvl 2 -> 3)
#
# Note: PowerPoint's TextRange.IndentLevel is 1-based (level 0 in the OOXML
# <a:pPr lvl="n"/> attribute corresponds to IndentLevel = 1), so bumping the
# XML lvl by one means bumping IndentLevel by one as well.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)
$sh = $s.Shapes.Item(2)
$tr = $sh.TextFrame.TextRange

$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $para = $tr.Paragraphs($i, 1)
    $txt = $para.Text
    if ($txt -like "Data architectures and methodologies for Digital Twins*") {
        $para.IndentLevel = 3
    } elseif ($txt -like "Precision Agriculture:*") {
        $para.IndentLevel = 3
    } elseif ($txt -like "Irrigation optimization*") {
        $para.IndentLevel = 4
    } elseif ($txt -like "Data Platform for italian agriculture domain*") {
        $para.IndentLevel = 4
    }
}
